$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(1.02,1.02579021841775,1.030664174593956,1.029418373908761,1.035993916326789,1.032230558027885,1.030957273845432,1.033474531516214,1.032232342361987,1.038788929057999,1.032421351426751),
    @(1.02,1.026938987174463,1.031537035326647,1.030514945218111,1.037337222847773,1.032509896211723,1.031744499328278,1.034155621235942,1.033136278540924,1.039940328501029,1.033209694860072),
    @(1.02,1.027681908650085,1.03210111613187,1.03122446411666,1.038206035137431,1.032688645428604,1.032252977048416,1.034594980859085,1.03372056594895,1.040684433770066,1.033718894676896),
    @(1.02,1.027994136838721,1.032338084304437,1.031522738447562,1.038571191652642,1.03276331280978,1.0324665242809,1.034779364907675,1.033966053100153,1.040997036034337,1.033932745170945),
    @(1.02,1.028046555749647,1.032377862269743,1.03157281962074,1.038632497715586,1.032775821724276,1.032502367116742,1.034810304895109,1.034007262868371,1.041049510521913,1.033968638907725),
    @(1.02,1.027686081033336,1.032104283184941,1.031228449700491,1.038210914735311,1.032689645018743,1.032255831327627,1.034597445875069,1.033723846737624,1.040688611636726,1.033721753009511),
    @(1.02,1.026178535942843,1.030959311092136,1.029788973521096,1.036447976936339,1.032325376244124,1.03122350937452,1.033704989120018,1.032537960808455,1.039178243297637,1.032687965040822),
    @(1.02,1.023518835084778,1.028936193485074,1.027252100978465,1.033338288238755,1.031668145933797,1.029397412285581,1.032121987812693,1.030443482112996,1.036509573380425,1.03085927468456),
    @(1.02,1.021743420999864,1.027583687451162,1.025560549148522,1.031262851913147,1.031219654370039,1.028175231714553,1.031059628510025,1.029043866052592,1.034725461581438,1.029635358476935),
    @(1.02,1.020974078853637,1.026997134819227,1.024827997686666,1.030363575886646,1.031022994823207,1.027644864019995,1.030597937538422,1.028437017911592,1.033951701149609,1.029104237599417),
    @(1.02,1.020688221901199,1.026779125521276,1.024555878947154,1.03002945117504,1.030949576583877,1.027447686801852,1.030426191102346,1.028211484963535,1.033664104439715,1.028906780367008),
    @(1.02,1.020749543256147,1.02682589549611,1.024614250115563,1.030101126352293,1.030965341803867,1.02748998992888,1.030463042800936,1.028259868116392,1.033725803451089,1.028949143569329),
    @(1.02,1.020950451659157,1.026979116914305,1.02480550461363,1.030335958963439,1.031016933602508,1.027628568868857,1.030583746109321,1.028418377796382,1.033927932161536,1.029087919307295),
    @(1.02,1.021074226086468,1.027073503431384,1.024923340582344,1.030480634579104,1.03104867194946,1.027713928689175,1.030658081730497,1.028516024558504,1.034052445374712,1.029173400348348),
    @(1.02,1.021794467355606,1.027622595787107,1.025609163956563,1.031322520982393,1.031232654143327,1.028210406004395,1.031090233880346,1.029084123468289,1.034776787370679,1.029670582718303),
    @(1.02,1.022246099691352,1.027966782937698,1.026039335288353,1.031850451300359,1.031347402364472,1.028521522743431,1.031360860317427,1.029440260243958,1.035230816946093,1.029982141278794),
    @(1.02,1.022509473769626,1.028167453883081,1.026290237575182,1.032158326631843,1.031414095770531,1.028702880488078,1.031518549888011,1.029647911220021,1.035495526416664,1.030163756572241),
    @(1.02,1.022599268198421,1.028235862662847,1.02637578723377,1.032263294522992,1.031436796267491,1.028764699943306,1.031572290468914,1.029718701725183,1.035585765533975,1.030225663818186),
    @(1.02,1.022197649565983,1.027929863981482,1.025993182957281,1.03179381534601,1.031335115514364,1.028488154393065,1.031331841452017,1.029402058173849,1.035182116127782,1.029948725541544),
    @(1.02,1.020891291619363,1.026934000830534,1.024749185400568,1.030266809224503,1.031001751325718,1.027587765668729,1.030548209009864,1.028371704041379,1.033868415498851,1.029047058161942),
    @(1.02,1.020069416231304,1.026307065892047,1.023966938375534,1.029306176766813,1.03079001029365,1.027020642078046,1.030054038928984,1.027723169808494,1.033041352632726,1.028479129190711),
    @(1.02,1.020505157691265,1.026639491787465,1.024381632042539,1.029815478912935,1.030902461428454,1.027321381575574,1.030316147476488,1.028067037872343,1.033479898395928,1.028780295772821),
    @(1.02,1.02221954226272,1.02794654633893,1.026014037237547,1.031819406872628,1.031340668146016,1.028503232458003,1.031344954329157,1.02941932028419,1.035204122295988,1.029963825019063),
    @(1.02,1.024206824614848,1.029459877602612,1.02790799088313,1.034142611214777,1.031839876167133,1.029870340172821,1.03253246631117,1.030985530829961,1.037200358326287,1.031332874183646)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $rowVals[$j]
    }
}
